# Fix mojibake "Â±" (double UTF-8 encoded ±) back to the proper "±"
# character in the benchmark result columns (f1_score_weighted,
# training_time, test_time) for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mojibake = "Â±"
$fixed = "±"

$range = $ws.Range("B2:D17")
foreach ($cell in $range.Cells) {
    $value = $cell.Value2
    if ($null -ne $value -and $value.Contains($mojibake)) {
        $cell.Value = $value.Replace($mojibake, $fixed)
    }
}
